# Update forecast figures after removing Auto ARIMA model from the pipeline.
# Updates the "Forecast Comparison" sheet (Prophet/Amazon forecast columns)
# and the "Summary" sheet (aggregate totals that depend on those columns).

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# row -> C (Prophet Forecast), D (Amazon Mean Forecast),
#        E (Amazon P70 Forecast), F (Amazon P80 Forecast), G (Amazon P90 Forecast)
$data = @{
    2  = @(594, 92,  110, 128, 156)
    3  = @(466, 100, 121, 144, 181)
    4  = @(464, 105, 127, 151, 188)
    5  = @(380, 105, 127, 150, 188)
    6  = @(209, 106, 129, 155, 197)
    7  = @(101, 105, 126, 151, 189)
    8  = @(140, 108, 132, 160, 205)
    9  = @(246, 109, 133, 162, 209)
    10 = @(293, 107, 130, 156, 198)
    11 = @(263, 108, 131, 160, 205)
    12 = @(235, 110, 134, 163, 210)
    13 = @(266, 117, 143, 175, 226)
    14 = @(312, 115, 141, 172, 221)
    15 = @(311, 113, 138, 170, 222)
    16 = @(271, 111, 136, 167, 218)
    17 = @(260, 110, 134, 166, 217)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $wsForecast.Cells.Item($row, 3).Value = $values[0]  # C - Prophet Forecast
    $wsForecast.Cells.Item($row, 4).Value = $values[1]  # D - Amazon Mean Forecast
    $wsForecast.Cells.Item($row, 5).Value = $values[2]  # E - Amazon P70 Forecast
    $wsForecast.Cells.Item($row, 6).Value = $values[3]  # F - Amazon P80 Forecast
    $wsForecast.Cells.Item($row, 7).Value = $values[4]  # G - Amazon P90 Forecast
}

# Update the summary totals that reflect the recomputed forecast.
# These cells hold text values (not numbers), so use a leading apostrophe
# to force text storage, matching the original cell type.
$wsSummary.Range("B10").Value = "'2600"
$wsSummary.Range("B11").Value = "'1904"
$wsSummary.Range("B12").Value = "'594"
